# Updated cryptos list (Price / Volume(1h) columns) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.216.78"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "2.280.16"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'112.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'265.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "'47.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'0.0932"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "'9.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.98%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").Value = "'15.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "2.611.21"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "'0.865"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").Value = "2.274.31"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "43.196.13"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "'0.0000108"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "'6.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.87%  "
$ws.Range("D21").Value = "'71.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "'2.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'234.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "'9.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("D25").Value = "'2.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'3.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'40.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "'173.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").Value = "'21.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "'0.0904"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "'5.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("D37").Value = "'4.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("D38").Value = "'0.0368"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.44%  "
$ws.Range("D39").Value = "'3.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.44%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("D41").Value = "'2.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.69%  "
$ws.Range("D42").Value = "'76.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.41%  "
$ws.Range("D43").Value = "'14.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").Value = "'6.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.96%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'1.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").Value = "'8.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "'103.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "'0.0996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.55%  "
